$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "first" column (column F) entirely - this shifts the
# "pts_per_game" column (G) left into F, matching the diff where the
# sheet dimension shrinks from A1:G6 to A1:F6.
$ws.Columns("F").Delete()
